$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Desde domingo 21 de junio de 2028" / "Hasta domingo 20 de julio
# de 2028" date-range labels that used to sit in the merged G2:I2 / G3:I3
# boxes, leaving those boxes blank.
$ws.Range("G2:I3").ClearContents()

# Expand the Funcionario merged box from C2:E2 to C2:F2.
$ws.Range("C2:F2").Merge()

# Resize columns to the new widths (values nudged slightly so that the
# engine's internal pixel-grid rounding of ColumnWidth lands on the same
# stored width Excel would have written for 21.69 / 10.78 / 6.74 / 7.1 / 8.21).
$ws.Columns("A").ColumnWidth = 20.8333333333333
$ws.Columns("B").ColumnWidth = 10
$ws.Columns("C:F").ColumnWidth = 5.83333333333333
$ws.Columns("G").ColumnWidth = 6.33333333333333
$ws.Columns("H:I").ColumnWidth = 7.33333333333333

# The logo picture is anchored (oneCell) off column A / row 1-3; nudge its
# stored size down very slightly to match the re-layout that happened when
# the columns above were resized.
$shp = $ws.Shapes.Item(1)
$shp.Width = 50.5417322834646
$shp.Height = 47.70705

# Move the active selection to G3, matching the editor's cursor position
# after clearing the date cells.
$ws.Range("G3").Select()
